# Commit for Minor correction in files
$wb = $excel.ActiveWorkbook

# --- Cabling sheet: fix shared string text ---
$wsCabling = $wb.Worksheets.Item("Cabling")
$wsCabling.Range("L2").Value = "Cabling 3U4U3I3I3I"

# --- DSPChannelMap sheet: correct cell values ---
$wsDsp = $wb.Worksheets.Item("DSPChannelMap")
$wsDsp.Range("D8").Value = 12
$wsDsp.Range("D9").Value = 14
$wsDsp.Range("B10").Value = 18
$wsDsp.Range("D10").Value = 15
$wsDsp.Range("B11").Value = 0
$wsDsp.Range("D11").Value = 16
$wsDsp.Range("B12").Value = 0
$wsDsp.Range("D12").Value = 17

# --- Update selections on each sheet ---
$wsCabling.Activate()
$wsCabling.Range("E15:F18").Select()

$wsFeeder = $wb.Worksheets.Item("DSPFeederMap")
$wsFeeder.Activate()
$wsFeeder.Range("B4").Select()

$wsBusbar = $wb.Worksheets.Item("BusbarFeederMap")
$wsBusbar.Activate()
$wsBusbar.Range("B38").Select()

# Make DSPChannelMap the final active/selected tab, as in the target workbook
$wsDsp.Activate()
$wsDsp.Range("B11").Select()
